$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backlog item: "Área de histórico" / "Área de histórico de calculo"
$ws.Cells.Item(9, 1).Value = "Área de histórico"
# B9 already carries a quote-prefixed style (matches B8's style); a leading
# apostrophe keeps that style (quotePrefix) instead of the engine silently
# switching to the "no quote prefix" variant of the same format.
$ws.Cells.Item(9, 2).Value = "'Área de histórico de calculo"

# A stray note in J11 (default styling)
$ws.Cells.Item(11, 10).Value = " "

# Zoom in on the sheet and leave the cursor parked on the new note
$excel.ActiveWindow.Zoom = 150
[void]$ws.Range("J11").Select()
